# Added Joins and Update, Delete, and Alter Table

$wb = $excel.ActiveWorkbook

# --- Sheet 1: FavFoodTable ---
$ws1 = $wb.Worksheets.Item("FavFoodTable")

# Change C2 formula from 10^10 to 10^9 (value drops from 10,000,000,000 to 1,000,000,000)
$ws1.Range("C2").Formula = "=10^9"

# Extend the table with a new (still empty) row 9 cell in column C, matching the
# number formatting already used for the rest of column C (style index 3).
$ws1.Range("C9").NumberFormat = $ws1.Range("C8").NumberFormat

# Move the active selection on sheet 1 to C8
$ws1.Activate()
$ws1.Range("C8").Select()

# --- Sheet 2: People ---
# Reset the lingering selection (previously H11) back to the sheet's default.
$ws2 = $wb.Worksheets.Item("People")
$ws2.Activate()
$ws2.Range("A1").Select()

# Keep sheet 1 as the active/visible tab, matching the saved workbook state.
$ws1.Activate()
